$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header row: rename existing "Swate template" -> "Swate template: Chromatography"
#    and add a new "Swate template: Mass spectrometry" header in column F.
$ws.Range("E1").Value = "Swate template: Chromatography"
$ws.Range("F1").Value = "Swate template: Mass spectrometry"
# Give F1 the same formatting (bold + green fill) as the other header cells.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 2. New note row under the header.
$ws.Range("A2").Value = "last updated: 23.4.24"

# 3. New extraction-parameter rows for the chromatography column (E).
$ws.Range("E6").Value = "solvent extraction"
$ws.Range("E7").Value = "derivatisation"

# 4. Update protocol type description to mention protocol ref as well.
$ws.Range("E18").Value = "protocol type, protocol ref"

# 5. Move the mass-spectrometry related values out of column E into the
#    newly added column F, leaving column E free for chromatography data.
$ws.Range("E26:E33").Cut($ws.Range("F26:F33"))
$ws.Range("E26:E33").Clear()

# 6. Unify the formatting of the "column model" / "column type" /
#    "(labelled extract name)" rows onto the same style used elsewhere
#    in the chromatography protocol block.
$ws.Range("C19").Copy()
$ws.Range("C20:D22").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
# re-apply values since PasteSpecial(formats) does not touch them, so they
# are already correct -- but make sure content stayed intact.
$ws.Range("C20").Value = "column model"
$ws.Range("D20").Value = "column model"
$ws.Range("C21").Value = "column type"
$ws.Range("D21").Value = "column type"
$ws.Range("C22").Value = "(labelled extract name)"
$ws.Range("D22").Value = "(labelled extract name)"

# 7. Update the column widths for E:F to better fit the new, longer text.
$ws.Range("E1:F1").ColumnWidth = 29.16

# 8. Restore selection near the top of the sheet (matches author's saved view).
$ws.Range("A3").Select()
